$rowsData = @(
    @{ Row=45; Cells=[ordered]@{ A=@{ T='text'; V='21CRB01268' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='POSSESSION DRUG PARAPHERNALIA' }; D=@{ T='text'; V='2925.14(C)' }; E=@{ T='text'; V='M4' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=46; Cells=[ordered]@{ A=@{ T='text'; V='21TRC09438' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='OVI ALCOHOL / DRUGS 1ST' }; D=@{ T='text'; V='4511.19A1A*' }; E=@{ T='text'; V='M1' }; F=@{ T='text'; V='Guilty' }; G=@{ T='text'; V='Guilty' }; H=@{ T='text'; V='30' }; I=@{ T='text'; V='25' }; J=@{ T='text'; V='10' }; K=@{ T='text'; V='5' } } },
    @{ Row=47; Cells=[ordered]@{ A=@{ T='text'; V='21TRC09438' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS' }; D=@{ T='text'; V='4510.11' }; E=@{ T='text'; V='M1' }; F=@{ T='text'; V='Guilty' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' }; J=@{ T='text'; V='None' }; K=@{ T='text'; V='None' } } },
    @{ Row=48; Cells=[ordered]@{ A=@{ T='text'; V='21TRC09438' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='REASONABLE CONTROL' }; D=@{ T='text'; V='4511.202' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='Guilty' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' }; J=@{ T='text'; V='None' }; K=@{ T='text'; V='None' } } },
    @{ Row=49; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS UCM' }; D=@{ T='text'; V='4510.111' }; E=@{ T='text'; V='UCM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=50; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='TAIL LIGHTS-REAR LICENSE PLATE' }; D=@{ T='text'; V='4513.05' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=51; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS UCM' }; D=@{ T='text'; V='4510.111' }; E=@{ T='text'; V='UCM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=52; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='TAIL LIGHTS-REAR LICENSE PLATE' }; D=@{ T='text'; V='4513.05' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=53; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS UCM' }; D=@{ T='text'; V='4510.111' }; E=@{ T='text'; V='UCM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=54; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='TAIL LIGHTS-REAR LICENSE PLATE' }; D=@{ T='text'; V='4513.05' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=55; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS UCM' }; D=@{ T='text'; V='4510.111' }; E=@{ T='text'; V='UCM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=56; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='TAIL LIGHTS-REAR LICENSE PLATE' }; D=@{ T='text'; V='4513.05' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=57; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS UCM' }; D=@{ T='text'; V='4510.111' }; E=@{ T='text'; V='UCM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=58; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='TAIL LIGHTS-REAR LICENSE PLATE' }; D=@{ T='text'; V='4513.05' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=59; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='DUS UCM' }; D=@{ T='text'; V='4510.111' }; E=@{ T='text'; V='UCM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } },
    @{ Row=60; Cells=[ordered]@{ A=@{ T='text'; V='21TRD09386' }; B=@{ T='text'; V='Bunner' }; C=@{ T='text'; V='TAIL LIGHTS-REAR LICENSE PLATE' }; D=@{ T='text'; V='4513.05' }; E=@{ T='text'; V='MM' }; F=@{ T='text'; V='No Contest' }; G=@{ T='text'; V='Guilty' }; H=@{ T='num'; V=0 }; I=@{ T='text'; V='0' } } }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($rowEntry in $rowsData) {
    $r = $rowEntry.Row
    foreach ($colLetter in $rowEntry.Cells.Keys) {
        $cellInfo = $rowEntry.Cells[$colLetter]
        $cell = $ws.Range("$colLetter$r")
        if ($cellInfo.T -eq 'num') {
            $cell.Value = $cellInfo.V
        }
        else {
            $v = [string]$cellInfo.V
            if ($v -match '^-?\d+(\.\d+)?$') {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $v
        }
    }
}

Write-Host "Added rows 45-60 to Sheet1"
